$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-12-30 Monday" "2024-12-31 Tuesday"
Replace-Text "13×12=156" "12×89=1068"
Replace-Text "39×67=2613" "66×55=3630"
Replace-Text "98×46=4508" "50×54=2700"
Replace-Text "62×79=4898" "58×99=5742"
Replace-Text "17×48=816" "11×88=968"
Replace-Text "49×21=1029" "36×84=3024"
Replace-Text "20×51=1020" "52×11=572"
Replace-Text "59×49=2891" "93×19=1767"
Replace-Text "36×64=2304" "86×23=1978"
Replace-Text "96×89=8544" "94×13=1222"
Replace-Text "77×29=2233" "58×71=4118"
Replace-Text "68×86=5848" "72×66=4752"
Replace-Text "99×89=8811" "59×74=4366"
Replace-Text "94×22=2068" "80×26=2080"
Replace-Text "32×93=2976" "73×84=6132"
Replace-Text "88×16=1408" "80×57=4560"
Replace-Text "23×82=1886" "49×54=2646"
Replace-Text "26×49=1274" "43×96=4128"
Replace-Text "23×95=2185" "90×53=4770"
Replace-Text "19×72=1368" "60×54=3240"
Replace-Text "34×54=1836" "17×73=1241"
Replace-Text "52×80=4160" "61×56=3416"
Replace-Text "36×23=828" "81×38=3078"
Replace-Text "14×24=336" "58×60=3480"
Replace-Text "59×82=4838" "35×84=2940"
